$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "74.873.87"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").Value = "2.817.91"
$ws.Range("E3").Value = "  +7.60%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'188.15"
$ws.Range("E5").Value = "  +1.85%  "
$ws.Range("D6").Value = "'596.59"
$ws.Range("E6").Value = "  +1.55%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  +3.47%  "
$ws.Range("E9").Value = "  -4.70%  "
$ws.Range("D10").Value = "2.813.01"
$ws.Range("E10").Value = "  +7.29%  "
$ws.Range("E11").Value = "  -1.19%  "
$ws.Range("E12").Value = "  +1.79%  "
$ws.Range("D13").Value = "'4.82"
$ws.Range("E13").Value = "  +0.36%  "
$ws.Range("D14").Value = "3.335.33"
$ws.Range("E14").Value = "  +7.15%  "
$ws.Range("D15").Value = "74.837.25"
$ws.Range("E15").Value = "  +1.18%  "
$ws.Range("D16").Value = "'27.00"
$ws.Range("E16").Value = "  +2.88%  "
$ws.Range("E17").Value = "  -1.82%  "
$ws.Range("D18").Value = "2.815.65"
$ws.Range("E18").Value = "  +6.21%  "
$ws.Range("D19").Value = "'8.94"
$ws.Range("E19").Value = "  -0.35%  "
$ws.Range("D20").Value = "'12.31"
$ws.Range("E20").Value = "  +3.77%  "
$ws.Range("D21").Value = "'374.82"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("E22").Value = "  -0.64%  "
$ws.Range("E23").Value = "  -0.31%  "
$ws.Range("D24").Value = "'6.17"
$ws.Range("E24").Value = "  -0.81%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").Value = "'70.68"
$ws.Range("E26").Value = "  +0.97%  "
$ws.Range("D27").Value = "2.959.48"
$ws.Range("E27").Value = "  +7.28%  "
$ws.Range("E28").Value = "  -0.89%  "
$ws.Range("D29").Value = "'9.54"
$ws.Range("E29").Value = "  +1.71%  "
$ws.Range("D30").Value = "'0.0000103"
$ws.Range("E30").Value = "  +7.77%  "
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  +0.69%  "
$ws.Range("D32").Value = "'513.25"
$ws.Range("E32").Value = "  +0.70%  "
$ws.Range("D33").Value = "'1.39"
$ws.Range("E33").Value = "  -0.66%  "
$ws.Range("D34").Value = "'7.87"
$ws.Range("E34").Value = "  -1.34%  "
$ws.Range("E35").Value = "  +2.51%  "
$ws.Range("E36").Value = "  -0.33%  "
$ws.Range("D37").Value = "'162.94"
$ws.Range("E37").Value = "  +1.91%  "
$ws.Range("D38").Value = "'20.08"
$ws.Range("E38").Value = "  +4.25%  "
$ws.Range("D39").Value = "'0.118"
$ws.Range("E39").Value = "  -4.14%  "
$ws.Range("D40").Value = "'19.33"
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").Value = "'182.29"
$ws.Range("E41").Value = "  +15.13%  "
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").Value = "'5.03"
$ws.Range("E43").Value = "  +2.14%  "
$ws.Range("E44").Value = "  +3.85%  "
$ws.Range("E45").Value = "  -0.55%  "
$ws.Range("E46").Value = "  +2.40%  "
$ws.Range("D47").Value = "'39.79"
$ws.Range("E47").Value = "  +2.40%  "
$ws.Range("E48").Value = "  -1.92%  "
$ws.Range("E49").Value = "  -7.63%  "
$ws.Range("D50").Value = "'0.567"
$ws.Range("E50").Value = "  +7.25%  "
$ws.Range("E51").Value = "  +2.66%  "
